$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.401.04'
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("D3").Value = '3.750.69'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '613.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").Value = '3.749.06'
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.526'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.40%  '
$ws.Range("E10").Value = '  -0.27%  '
$ws.Range("E11").Value = '  +3.51%  '
$ws.Range("E12").Value = '  -3.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '39.81'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.33%  '
$ws.Range("E14").Value = '  -0.30%  '
$ws.Range("D15").Value = '4.371.11'
$ws.Range("E15").Value = '  +0.00%  '
$ws.Range("D16").Value = '3.749.58'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").Value = '69.523.00'
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("E18").Value = '  -2.55%  '
$ws.Range("E19").Value = '  -2.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '499.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = '  -3.72%  '
$ws.Range("E23").Value = '  -0.99%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.56'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.72'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.87'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("E28").Value = '  +5.53%  '
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.48'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.04'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.91'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.40%  '
$ws.Range("E33").Value = '  -3.05%  '
$ws.Range("E34").Value = '  -2.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("E36").Value = '  +0.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.09'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.09%  '
$ws.Range("E38").Value = '  +2.54%  '
$ws.Range("E39").Value = '  +3.63%  '
$ws.Range("E40").Value = '  +12.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '445.06'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.33%  '
$ws.Range("E42").Value = '  -5.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.69'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.88%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '44.43'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.55'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.01%  '
$ws.Range("D46").Value = '2.946.82'
$ws.Range("E46").Value = '  -4.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0359'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.48%  '
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '137.93'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '27.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.79%  '
$ws.Range("E51").Value = '  -2.10%  '
